# Append a new row (row 65) of sensor data to each of the four worksheets.
$wb = $excel.ActiveWorkbook

$rowsData = @{
    "ROW35-FE-LIFTER"  = @("2025-03-07 00:42:06", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0x d", 400, "568631262647113770877196", 400, 13)
    "ROW35-MID-LIFTER" = @("2025-03-07 00:29:35", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,", "0x01,0x90,", "0x e", 400, "568631262647113770942732", 400, 14)
    "ROW02-FE-LIFTER"  = @("2025-03-07 00:51:45", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,", "0x01,0x90,", "0xff", 400, "568631262647113769959692", 400, 255)
    "ROW02-MID-LIFTER" = @("2025-03-07 00:41:15", "0x01,0x90 ", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x01,0x90,", "0x 3", 400, "568631262647113769959692", 400, 3)
}

foreach ($sheetName in $rowsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $rowsData[$sheetName]

    $ws.Cells.Item(65, 1).Value = $values[0]
    $ws.Cells.Item(65, 2).Value = $values[1]
    $ws.Cells.Item(65, 3).Value = $values[2]
    $ws.Cells.Item(65, 4).Value = $values[3]
    $ws.Cells.Item(65, 5).Value = $values[4]
    $ws.Cells.Item(65, 6).Value = $values[5]

    $gCell = $ws.Cells.Item(65, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $values[6]
    $gCell.Style = "Normal"

    $ws.Cells.Item(65, 8).Value = $values[7]
    $ws.Cells.Item(65, 9).Value = $values[8]
}
